$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text ("@") for the Price column cells that are being updated,
# so Excel keeps the value as a literal text string (matching the original inlineStr cells)
# instead of auto-converting single-dot decimal-looking values into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.146.47"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.47"
$ws.Range("E3").Value = "  -2.92%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.48"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("E7").Value = "  -3.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.05"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.264"
$ws.Range("E10").Value = "  -6.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08162"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.03"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.633"
$ws.Range("E14").Value = "  -6.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.440"
$ws.Range("E15").Value = "  -6.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001257"
$ws.Range("E16").Value = "  -4.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.602.10"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.49"
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06854"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("E20").Value = "  -7.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.577"
$ws.Range("E21").Value = "  -5.75%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.98"
$ws.Range("E23").Value = "  -5.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.144.11"
$ws.Range("E24").Value = "  -3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.816"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.08"
$ws.Range("E27").Value = "  -4.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.28"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.279"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.12"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.409"
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.920"
$ws.Range("E32").Value = "  -12.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.778.40"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07723"
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9443"
$ws.Range("E35").Value = "  -7.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02754"
$ws.Range("E36").Value = "  -5.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.263"
$ws.Range("E37").Value = "  -7.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2539"
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.15"
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08924"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.390"
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7117"
$ws.Range("E42").Value = "  -6.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.74"
$ws.Range("E43").Value = "  -4.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.60"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6628"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.303"
$ws.Range("E47").Value = "  -6.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.980"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.56"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07951"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("E51").Value = "  -0.71%  "
